$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Text)
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "29.380.86"
Set-TextValue $ws "E2" "  -0.04%  "
Set-TextValue $ws "D3" "1.847.42"
Set-TextValue $ws "E3" "  -0.06%  "
Set-TextValue $ws "D4" "0.9985"
Set-TextValue $ws "E4" "  -0.09%  "
Set-TextValue $ws "D5" "240.24"
Set-TextValue $ws "E5" "  -0.02%  "
Set-TextValue $ws "D6" "0.6278"
Set-TextValue $ws "E6" "  -0.31%  "
Set-TextValue $ws "D7" "0.9999"
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "D8" "0.07642"
Set-TextValue $ws "E8" "  +0.36%  "
Set-TextValue $ws "D9" "0.2906"
Set-TextValue $ws "E9" "  -0.98%  "
Set-TextValue $ws "D10" "24.71"
Set-TextValue $ws "E10" "  +0.80%  "
Set-TextValue $ws "D11" "0.07739"
Set-TextValue $ws "E11" "  -0.05%  "
Set-TextValue $ws "D12" "5.034"
Set-TextValue $ws "E13" "  -0.09%  "
Set-TextValue $ws "D14" "0.00001061"
Set-TextValue $ws "E14" "  -1.84%  "
Set-TextValue $ws "E15" "  -0.55%  "
Set-TextValue $ws "D16" "6.155"
Set-TextValue $ws "E16" "  -0.01%  "
Set-TextValue $ws "D17" "29.424.16"
Set-TextValue $ws "E17" "  +0.02%  "
Set-TextValue $ws "D18" "226.77"
Set-TextValue $ws "E18" "  -0.92%  "
Set-TextValue $ws "E19" "  -0.85%  "
Set-TextValue $ws "D21" "7.501"
Set-TextValue $ws "E21" "  +0.70%  "
Set-TextValue $ws "D22" "0.9994"
Set-TextValue $ws "E22" "  -0.11%  "
Set-TextValue $ws "D23" "158.28"
Set-TextValue $ws "E23" "  +0.61%  "
Set-TextValue $ws "D24" "0.1384"
Set-TextValue $ws "E24" "  -0.48%  "
Set-TextValue $ws "D25" "8.405"
Set-TextValue $ws "E25" "  +0.31%  "
Set-TextValue $ws "D26" "17.68"
Set-TextValue $ws "E26" "  +0.26%  "
Set-TextValue $ws "D27" "1.385"
Set-TextValue $ws "E27" "  +5.51%  "
Set-TextValue $ws "D28" "1.461"
Set-TextValue $ws "E28" "  -0.45%  "
Set-TextValue $ws "E29" "  -0.06%  "
Set-TextValue $ws "D31" "4.083"
Set-TextValue $ws "E31" "  +0.95%  "
Set-TextValue $ws "D32" "1.838"
Set-TextValue $ws "E32" "  -0.61%  "
Set-TextValue $ws "E33" "  +0.48%  "
Set-TextValue $ws "D34" "0.6945"
Set-TextValue $ws "E34" "  -2.05%  "
Set-TextValue $ws "D35" "2.579"
Set-TextValue $ws "E35" "  -0.17%  "
Set-TextValue $ws "E36" "  +0.28%  "
Set-TextValue $ws "D37" "1.229.33"
Set-TextValue $ws "E37" "  -0.31%  "
Set-TextValue $ws "D38" "2.714"
Set-TextValue $ws "E38" "  -2.10%  "
Set-TextValue $ws "D39" "6.389"
Set-TextValue $ws "E39" "  -1.31%  "
Set-TextValue $ws "D40" "0.9043"
Set-TextValue $ws "E40" "  -0.46%  "
Set-TextValue $ws "E42" "  +0.01%  "
Set-TextValue $ws "D43" "66.06"
Set-TextValue $ws "E43" "  +0.04%  "
Set-TextValue $ws "D44" "7.183"
Set-TextValue $ws "E44" "  +0.33%  "
Set-TextValue $ws "B45" "TheSandbox"
Set-TextValue $ws "C45" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D45" "0.4011"
Set-TextValue $ws "E45" "  -0.14%  "
Set-TextValue $ws "B46" "EnergySwap"
Set-TextValue $ws "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D46" "9.008"
Set-TextValue $ws "E46" "  +0.05%  "
Set-TextValue $ws "B47" "RenderToken"
Set-TextValue $ws "C47" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D47" "1.677"
Set-TextValue $ws "E47" "  -0.43%  "
Set-TextValue $ws "B48" "Algorand"
Set-TextValue $ws "C48" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D48" "0.1142"
Set-TextValue $ws "E48" "  +1.69%  "
Set-TextValue $ws "B49" "Cronos"
Set-TextValue $ws "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D49" "0.05701"
Set-TextValue $ws "E49" "  -0.12%  "
Set-TextValue $ws "B50" "Mantle"
Set-TextValue $ws "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D50" "0.4624"
Set-TextValue $ws "E50" "  +0.01%  "
Set-TextValue $ws "B51" "NEARProtocol"
Set-TextValue $ws "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D51" "1.343"
Set-TextValue $ws "E51" "  +0.21%  "
